# Shift the whole itinerary back by one day.
# Trip dates: C1 (departure) and C2 (return) each move one day earlier,
# and the day labels in column B (rows 7-11) follow suit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header dates: 20-24 May -> 19-23 May (serial dates, one day earlier each)
$ws.Range("C1").Value = 45796
$ws.Range("C2").Value = 45800

# Day labels in the itinerary table
$ws.Range("B7").Value  = "chiều 19"
$ws.Range("B8").Value  = "Ngày 20"
$ws.Range("B9").Value  = "Ngày 21"
$ws.Range("B10").Value = "Ngày 22"
$ws.Range("B11").Value = "Ngày 23"

# Restore the view state (scrolled/selected cell) left by the author on save
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("B12").Select()
